$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 51
$ws.Cells.Item(51, 1).Value = 50
$ws.Cells.Item(51, 2).Value = "Monday, Jan 09"
$ws.Cells.Item(51, 3).Value = "2:25 PM"
$ws.Cells.Item(51, 4).Value = "SK1756"
$ws.Cells.Item(51, 5).Value = "Copenhagen"
$ws.Cells.Item(51, 6).Value = "(CPH)"
$ws.Cells.Item(51, 7).Value = "SAS "
$ws.Cells.Item(51, 8).Value = "CRJ9"
$ws.Cells.Item(51, 9).Value = "(EI-FPV)"
$ws.Cells.Item(51, 10).Value = "2:46 PM"
$ws.Cells.Item(51, 12).Value = "0 hours, 21 minutes"
$ws.Cells.Item(51, 11).Borders.LineStyle = -4142
$ws.Cells.Item(51, 13).Borders.LineStyle = -4142

# Row 52
$ws.Cells.Item(52, 1).Value = 51
$ws.Cells.Item(52, 2).Value = "Monday, Jan 09"
$ws.Cells.Item(52, 3).Value = "2:45 PM"
$ws.Cells.Item(52, 4).Value = "LO3944"
$ws.Cells.Item(52, 5).Value = "Warsaw"
$ws.Cells.Item(52, 6).Value = "(WAW)"
$ws.Cells.Item(52, 7).Value = "LOT "
$ws.Cells.Item(52, 8).Value = "E170"
$ws.Cells.Item(52, 9).Value = "(SP-LDG)"
$ws.Cells.Item(52, 10).Value = "2:50 PM"
$ws.Cells.Item(52, 12).Value = "0 hours, 5 minutes"
$ws.Cells.Item(52, 11).Borders.LineStyle = -4142
$ws.Cells.Item(52, 13).Borders.LineStyle = -4142

# Row 53
$ws.Cells.Item(53, 1).Value = 52
$ws.Cells.Item(53, 2).Value = "Monday, Jan 09"
$ws.Cells.Item(53, 3).Value = "3:15 PM"
$ws.Cells.Item(53, 4).Value = "KL1274"
$ws.Cells.Item(53, 5).Value = "Amsterdam"
$ws.Cells.Item(53, 6).Value = "(AMS)"
$ws.Cells.Item(53, 7).Value = "KLM "
$ws.Cells.Item(53, 8).Value = "E75L"
$ws.Cells.Item(53, 9).Value = "(PH-EXW)"
$ws.Cells.Item(53, 10).Value = "3:11 PM"
$ws.Cells.Item(53, 12).Value = "0 hours, -4 minutes"
$ws.Cells.Item(53, 11).Borders.LineStyle = -4142
$ws.Cells.Item(53, 13).Borders.LineStyle = -4142

# Row 54
$ws.Cells.Item(54, 1).Value = 53
$ws.Cells.Item(54, 2).Value = "Monday, Jan 09"
$ws.Cells.Item(54, 3).Value = "3:55 PM"
$ws.Cells.Item(54, 4).Value = "FR7948"
$ws.Cells.Item(54, 5).Value = "Bristol"
$ws.Cells.Item(54, 6).Value = "(BRS)"
$ws.Cells.Item(54, 7).Value = "Ryanair "
$ws.Cells.Item(54, 8).Value = "B738"
$ws.Cells.Item(54, 9).Value = "(SP-RSX)"
$ws.Cells.Item(54, 10).Value = "4:00 PM"
$ws.Cells.Item(54, 12).Value = "0 hours, 5 minutes"
$ws.Cells.Item(54, 11).Borders.LineStyle = -4142
$ws.Cells.Item(54, 13).Borders.LineStyle = -4142

# Row 55
$ws.Cells.Item(55, 1).Value = 54
$ws.Cells.Item(55, 2).Value = "Monday, Jan 09"
$ws.Cells.Item(55, 3).Value = "4:30 PM"
$ws.Cells.Item(55, 4).Value = "W91901"
$ws.Cells.Item(55, 5).Value = "London"
$ws.Cells.Item(55, 6).Value = "(LTN)"
$ws.Cells.Item(55, 7).Value = "Wizz Air "
$ws.Cells.Item(55, 8).Value = "A320"
$ws.Cells.Item(55, 9).Value = "(G-WUKD)"
$ws.Cells.Item(55, 10).Value = "4:34 PM"
$ws.Cells.Item(55, 12).Value = "0 hours, 4 minutes"
$ws.Cells.Item(55, 11).Borders.LineStyle = -4142
$ws.Cells.Item(55, 13).Borders.LineStyle = -4142
